# Update the micro_reg regression table: the "rincvar I..IIII" column group
# is renamed to "inciqr I..IIII" and its coefficients/SEs/N/R2 are refreshed
# to a new regression run's output.
#
# All of the touched cells hold plain (non-formula) text in the source file,
# including numeric-looking strings such as "57428" or "-0.01" that must stay
# text (shared-string) rather than become real numbers. Typing a numeric
# string straight into Range.Value (the normal COM path) makes Excel coerce
# it to a number, so instead we stage each value in a scratch cell that is
# explicitly number-formatted as Text ("@"), copy it, and paste-special just
# the values into the destination. That preserves the destination cell's own
# existing style (most of these cells carry no explicit style at all) while
# still landing the value as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell well outside the table's used range (A1:I17), fully deleted
# again after every use so it never ends up holding a leftover value/style
# on save.
$scratchAddr = "K1"

function Set-TextValue {
    param([string]$addr, [string]$val)
    $scratch = $ws.Range($scratchAddr)
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    # xlShiftToLeft: fully remove the scratch cell (not just its contents) so
    # no stray styled/empty cell - and no enlarged used-range/dimension -
    # survives the save.
    $scratch.Delete(-4159)
}

# Header row: rename "rincvar I..IIII" -> "inciqr I..IIII"
Set-TextValue "F1" "inciqr I"
Set-TextValue "G1" "inciqr II"
Set-TextValue "H1" "inciqr III"
Set-TextValue "I1" "inciqr IIII"

# parttime=yes coefficient row + SE row
Set-TextValue "F2" "0.03***"
Set-TextValue "G2" "0.04***"
Set-TextValue "H2" "-0.01"
Set-TextValue "I2" "0.02"

Set-TextValue "F3" "(0.01)"
Set-TextValue "G3" "(0.01)"
Set-TextValue "H3" "(0.01)"
Set-TextValue "I3" "(0.01)"

# selfemp=yes coefficient row + SE row
Set-TextValue "F4" "0.67***"
Set-TextValue "G4" "-0.00***"
Set-TextValue "I4" "-0.00***"

Set-TextValue "F5" "(0.01)"

# UEprobAgg coefficient row
Set-TextValue "G6" "0.00***"
Set-TextValue "H6" "0.00***"
Set-TextValue "I6" "0.00***"

# HHinc_gr=low inc coefficient + SE row
Set-TextValue "H10" "0.18***"
Set-TextValue "I10" "0.18***"

Set-TextValue "I11" "(0.01)"

# educ_gr=low educ coefficient + SE row
Set-TextValue "H12" "0.04***"
Set-TextValue "I12" "0.04***"

Set-TextValue "H13" "(0.01)"
Set-TextValue "I13" "(0.01)"

# gender=male coefficient row
Set-TextValue "I14" "-0.08***"

# N row
Set-TextValue "F16" "57428"
Set-TextValue "G16" "50772"
Set-TextValue "H16" "50772"
Set-TextValue "I16" "43094"

# R2 row
Set-TextValue "F17" "0.06"
Set-TextValue "H17" "0.02"
Set-TextValue "I17" "0.02"

Write-Output "micro_reg table updated"
